# Actualización automática de tasas-transfi.xlsx
# Updates the "Conversión del día" note on Hoja1 and the tasas (N10/O10/N12/O12)
# values on the "tasas" sheet.

$wb = $excel.ActiveWorkbook

# --- Hoja1: update the conversion note text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 6.52 = 25749.74 pesos`n✅ 25749.74 pesos = 6.49 = 964.81 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- tasas: update the rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 153.4
$wsTasas.Range("O10").Value = 3950.01
$wsTasas.Range("N12").Value = 3969.99
$wsTasas.Range("O12").Value = 148.751
